$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = "Daniel iba a Bolivia."
$ws.Range("E5").Value = "Emilio ama la marcha."
$ws.Range("E6").Value = "Mariano habla del tiempo."
$ws.Range("E8").Value = "María bebe el vino."
$ws.Range("E9").Value = "Marta abre el regalo."
$ws.Range("E10").Value = "Manuela vende el carro."
$ws.Range("E13").Value = "Emilio ama la marcha."
$ws.Range("D14").Value = "¿De qué habla Mariano?"
$ws.Range("E14").Value = "Mariano habla del tiempo."
$ws.Range("D16").Value = "¿Qué bebe María?"
$ws.Range("E16").Value = "María bebe el vino."
$ws.Range("D17").Value = "¿Qué abre Marta?"
$ws.Range("E17").Value = "Marta abre el regalo."
$ws.Range("D18").Value = "¿Qué vende Manuela?"
$ws.Range("E18").Value = "Manuela vende el carro."
$ws.Range("E21").Value = "¿Emilio ama la marcha?"
$ws.Range("E22").Value = "¿Mariano habla del tiempo?"
$ws.Range("E24").Value = "¿María bebe el vino?"
$ws.Range("E25").Value = "¿Marta abre el regalo?"
$ws.Range("E26").Value = "¿Manuela vende el carro?"
$ws.Range("E36").Value = "La niña lava el plato."
$ws.Range("E37").Value = "Mi madre come la fruta."
$ws.Range("E38").Value = "El niño oye el río."
$ws.Range("E39").Value = "Mi tía odia la lluvia."
$ws.Range("E41").Value = "La maestra vive en Paris."
$ws.Range("D44").Value = "¿Qué lava la niña?"
$ws.Range("E44").Value = "La niña lava el plato."
$ws.Range("D45").Value = "¿Qué come tu madre?"
$ws.Range("E45").Value = "Mi madre come la fruta."
$ws.Range("D46").Value = "¿Qué oye el niño?"
$ws.Range("E46").Value = "El niño oye el río."
$ws.Range("D47").Value = "¿Qué odia tu tía?"
$ws.Range("E47").Value = "Mi tía odia la lluvia."
$ws.Range("D49").Value = "¿Dónde vive la maestra?"
$ws.Range("E49").Value = "La maestra vive en Paris."
$ws.Range("E52").Value = "¿La niña lava el plato?"
$ws.Range("E53").Value = "¿Mi madre come la fruta?"
$ws.Range("E54").Value = "¿El niño oye el río?"
$ws.Range("E55").Value = "¿Mi tía odia la lluvia?"
$ws.Range("E57").Value = "¿La maestra vive en Paris?"
$ws.Range("E59").Value = "¿Cuándo miraba la luna?"
$ws.Range("E61").Value = "¿Cuándo comía la fruta?"
$ws.Range("E63").Value = "¿Por qué odiaba la lluvia?"
$ws.Range("E64").Value = "¿Por qué desayuna muy bien?"
$ws.Range("E65").Value = "¿Por qué vivía en Paris?"
$ws.Range("E66").Value = "¿Por qué venía del lago?"
